{"js": "// The edit:\n//  - Paragraph 1's two split runs (\"...the pag\" / \"e is too big...\") are\n//    merged into a single run, and the `_GoBack` bookmark that used to sit\n//    between them is removed from this spot.\n//  - Two brand-new paragraphs are inserted right after paragraph 1:\n//      * \"My website page name will be Lite's profile.\"\n//      * A long paragraph that starts with new \"home content\" text, then\n//        continues with the (previously paragraph-2) \"background color\"\n//        text, further new sentences (incl. one run with an eastAsia font\n//        hint), the `_GoBack` bookmark re-inserted mid-word (\"so|me\"), and\n//        finally the unchanged tail (\"...are interested in me. ... bootstrap\n//        image for Instagram.\") followed by the second picture.\n//\n// Strategy: read the whole body's OOXML just to pull out the two\n// <w:drawing> blocks (so the pictures/relationships carry over untouched),\n// then rebuild the whole body as OOXML and swap it in with a single\n// Body.insertOoxml(..., \"Replace\") call.\n\nconst body = context.document.body;\n\n// Pull the raw OOXML for the *whole body* in one shot so both pictures'\n// <w:drawing> blocks (and their r:embed relationship ids, which getOoxml()\n// renumbers per-call) come from the same, mutually-consistent snapshot.\nconst bodyOoxml = body.getOoxml();\nawait context.sync();\n\nfunction extractDrawings(ooxml) {\n  const matches = ooxml.value.match(/<w:drawing>[\\s\\S]*?<\\/w:drawing>/g);\n  if (!matches || matches.length < 2) {\n    throw new Error(\"Could not locate both <w:drawing> blocks\");\n  }\n  return matches;\n}\n\nconst [drawing1, drawing2] = extractDrawings(bodyOoxml); // picture 1, picture 2\n\nconst RSID1 = \"_GoBack\";\n\n// Paragraph 1: picture + the single merged sentence run.\nconst newPara1 =\n  \"<w:p>\" +\n  \"<w:r><w:rPr><w:noProof/></w:rPr>\" + drawing1 + \"</w:r>\" +\n  \"<w:r><w:t>This is the first part of the assignment because the page is too big. The second part is on second page.</w:t></w:r>\" +\n  \"</w:p>\";\n\n// New paragraph 2: \"My website page name will be Lite's profile.\"\nconst newPara2 =\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">My website page name will be </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>Lite\\u2019s</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> profile.</w:t></w:r>' +\n  \"</w:p>\";\n\n// New paragraph 3: the big \"home content\" paragraph, re-using the original\n// \"background color ... Instagram.\" text plus new surrounding sentences,\n// with the _GoBack bookmark now sitting inside \"so|me\" and the second\n// picture at the very end.\nconst newPara3 =\n  \"<w:p>\" +\n  \"<w:r><w:t>For the home content I am thinking about putting images about myself. Probably a gallery or even just a simple image following some quotes or just a small introduction to the website. For the header I will put navigations to each title in the page. Each section would have a topic and I haven\\u2019t decided if I want a background image for each section or simply a background color. If I were to use a</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> background color I think I will use </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/>' +\n  \"<w:r><w:t>rgb</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  \"<w:r><w:t>(</w:t></w:r>\" +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  \"<w:r><w:t>146,172,209)</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> or </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>Morandi</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> colors. I will try with different colors for the background color and the titles of course. For the education and Club experience part I will add information about the recent colleges I went and so</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"' + RSID1 + '\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  \"<w:r><w:t>me important courses worth pointing out. I will also include the experiences I had in the club since I was captain for the softball team in college. For the contact me form, I will include forms with First and Last Name, email and phone number.</w:t></w:r>\" +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\"> I will also add a text column if they want to send me some messages.</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> It will be a simple form for them to contact me if they </w:t></w:r>' +\n  \"<w:r><w:lastRenderedPageBreak/><w:t>are interested in me. For the footer part, I will also include a bootstrap image for Instagram.</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:noProof/></w:rPr>\" + drawing2 + \"</w:r>\" +\n  \"</w:p>\";\n\nconst newBodyXml = newPara1 + newPara2 + newPara3;\n\nconst nsDecl =\n  'xmlns:wpc=\"http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas\" ' +\n  'xmlns:cx=\"http://schemas.microsoft.com/office/drawing/2014/chartex\" ' +\n  'xmlns:mc=\"http://schemas.openxmlformats.org/markup-compatibility/2006\" ' +\n  'xmlns:o=\"urn:schemas-microsoft-com:office:office\" ' +\n  'xmlns:r=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\" ' +\n  'xmlns:m=\"http://schemas.openxmlformats.org/officeDocument/2006/math\" ' +\n  'xmlns:v=\"urn:schemas-microsoft-com:vml\" ' +\n  'xmlns:wp14=\"http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing\" ' +\n  'xmlns:wp=\"http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing\" ' +\n  'xmlns:w10=\"urn:schemas-microsoft-com:office:word\" ' +\n  'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" ' +\n  'xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\" ' +\n  'xmlns:w15=\"http://schemas.microsoft.com/office/word/2012/wordml\" ' +\n  'xmlns:w16se=\"http://schemas.microsoft.com/office/word/2015/wordml/symex\" ' +\n  'xmlns:wpg=\"http://schemas.microsoft.com/office/word/2010/wordprocessingGroup\" ' +\n  'xmlns:wpi=\"http://schemas.microsoft.com/office/word/2010/wordprocessingInk\" ' +\n  'xmlns:wne=\"http://schemas.microsoft.com/office/word/2006/wordml\" ' +\n  'xmlns:wps=\"http://schemas.microsoft.com/office/word/2010/wordprocessingShape\" ' +\n  'mc:Ignorable=\"w14 w15 w16se wp14\"';\n\nconst pkg =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  \"<w:document \" + nsDecl + \">\" +\n  \"<w:body>\" + newBodyXml + \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nbody.insertOoxml(pkg, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The edit:\n#  - Paragraph 1's two split runs (\"...the pag\" / \"e is too big...\") are\n#    merged into a single run, and the `_GoBack` bookmark that used to sit\n#    between them is removed from this spot.\n#  - Two brand-new paragraphs are inserted right after paragraph 1:\n#      * \"My website page name will be Lite's profile.\"\n#      * A long paragraph that starts with new \"home content\" text, then\n#        continues with the (previously paragraph-2) \"background color\"\n#        text, further new sentences (incl. one run with an eastAsia font\n#        hint), the `_GoBack` bookmark re-inserted mid-word (\"so|me\"), and\n#        finally the unchanged tail (\"...are interested in me. ... bootstrap\n#        image for Instagram.\") followed by the second picture.\n#\n# Strategy: read the whole document's WordOpenXML in one shot to pull out\n# both pictures' <w:drawing> blocks (so the pictures/relationships carry\n# over untouched, and r:embed ids - which get renumbered per WordOpenXML\n# read - stay mutually consistent), then rebuild the whole body as OOXML\n# and swap it in with a single Range.InsertXML(...) call on $d.Content.\n\n$d = $word.ActiveDocument\n\n# Read the whole document's OOXML once, so both <w:drawing> blocks come\n# from the same, mutually-consistent relationship-id snapshot.\n$fullXml = $d.Content.WordOpenXML\n$drawingMatches = [regex]::Matches($fullXml, '<w:drawing>[\\s\\S]*?</w:drawing>')\nif ($drawingMatches.Count -lt 2) {\n    throw \"Could not locate both <w:drawing> blocks\"\n}\n$drawing1 = $drawingMatches[0].Value\n$drawing2 = $drawingMatches[1].Value\n\n# Paragraph 1: picture + the single merged sentence run.\n$newPara1 = '<w:p><w:r><w:rPr><w:noProof/></w:rPr>' + $drawing1 + '</w:r>' + `\n    '<w:r><w:t>This is the first part of the assignment because the page is too big. The second part is on second page.</w:t></w:r></w:p>'\n\n# New paragraph 2: \"My website page name will be Lite's profile.\"\n$newPara2 = @'\n<w:p><w:r><w:t xml:space=\"preserve\">My website page name will be </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Lite\u2019s</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> profile.</w:t></w:r></w:p>\n'@\n\n# New paragraph 3 (everything up to, but not including, the trailing run +\n# second picture, which we append afterwards since they need PowerShell\n# variables spliced in).\n$newPara3Head = @'\n<w:p><w:r><w:t>For the home content I am thinking about putting images about myself. Probably a gallery or even just a simple image following some quotes or just a small introduction to the website. For the header I will put navigations to each title in the page. Each section would have a topic and I haven\u2019t decided if I want a background image for each section or simply a background color. If I were to use a</w:t></w:r><w:r><w:t xml:space=\"preserve\"> background color I think I will use </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/><w:r><w:t>rgb</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>146,172,209)</w:t></w:r><w:r><w:t xml:space=\"preserve\"> or </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Morandi</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> colors. I will try with different colors for the background color and the titles of course. For the education and Club experience part I will add information about the recent colleges I went and so</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t>me important courses worth pointing out. I will also include the experiences I had in the club since I was captain for the softball team in college. For the contact me form, I will include forms with First and Last Name, email and phone number.</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\"> I will also add a text column if they want to send me some messages.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> It will be a simple form for them to contact me if they </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>are interested in me. For the footer part, I will also include a bootstrap image for Instagram.</w:t></w:r>\n'@\n\n$newPara3 = $newPara3Head + '<w:r><w:rPr><w:noProof/></w:rPr>' + $drawing2 + '</w:r></w:p>'\n\n$newBodyXml = $newPara1 + $newPara2 + $newPara3\n\n$nsDecl = 'xmlns:wpc=\"http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas\" ' + `\n    'xmlns:cx=\"http://schemas.microsoft.com/office/drawing/2014/chartex\" ' + `\n    'xmlns:mc=\"http://schemas.openxmlformats.org/markup-compatibility/2006\" ' + `\n    'xmlns:o=\"urn:schemas-microsoft-com:office:office\" ' + `\n    'xmlns:r=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\" ' + `\n    'xmlns:m=\"http://schemas.openxmlformats.org/officeDocument/2006/math\" ' + `\n    'xmlns:v=\"urn:schemas-microsoft-com:vml\" ' + `\n    'xmlns:wp14=\"http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing\" ' + `\n    'xmlns:wp=\"http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing\" ' + `\n    'xmlns:w10=\"urn:schemas-microsoft-com:office:word\" ' + `\n    'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" ' + `\n    'xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\" ' + `\n    'xmlns:w15=\"http://schemas.microsoft.com/office/word/2012/wordml\" ' + `\n    'xmlns:w16se=\"http://schemas.microsoft.com/office/word/2015/wordml/symex\" ' + `\n    'xmlns:wpg=\"http://schemas.microsoft.com/office/word/2010/wordprocessingGroup\" ' + `\n    'xmlns:wpi=\"http://schemas.microsoft.com/office/word/2010/wordprocessingInk\" ' + `\n    'xmlns:wne=\"http://schemas.microsoft.com/office/word/2006/wordml\" ' + `\n    'xmlns:wps=\"http://schemas.microsoft.com/office/word/2010/wordprocessingShape\" ' + `\n    'mc:Ignorable=\"w14 w15 w16se wp14\"'\n\n$pkg = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n    '<pkg:xmlData>' + `\n    '<w:document ' + $nsDecl + '>' + `\n    '<w:body>' + $newBodyXml + '</w:body>' + `\n    '</w:document>' + `\n    '</pkg:xmlData>' + `\n    '</pkg:part>' + `\n    '</pkg:package>'\n\n[void]$d.Content.InsertXML($pkg)\n"}
